$d = $word.ActiveDocument

$replacements = @(
    @('2025-09-13 Saturday', '2025-09-14 Sunday'),
    @('437×7=3059', '510×7=3570'),
    @('906×7=6342', '102×8=816'),
    @('599×3=1797', '757×8=6056'),
    @('768×7=5376', '505×4=2020'),
    @('262×2=524', '112×9=1008'),
    @('977×4=3908', '602×7=4214'),
    @('557×9=5013', '390×2=780'),
    @('604×7=4228', '998×3=2994'),
    @('173×6=1038', '336×8=2688'),
    @('978×3=2934', '199×5=995'),
    @('404×7=2828', '467×6=2802'),
    @('531×6=3186', '428×3=1284'),
    @('290×4=1160', '270×3=810'),
    @('882×2=1764', '512×4=2048'),
    @('974×5=4870', '787×2=1574'),
    @('207×6=1242', '799×5=3995'),
    @('203×9=1827', '271×8=2168'),
    @('201×3=603', '522×4=2088'),
    @('970×5=4850', '685×4=2740'),
    @('255×4=1020', '166×4=664'),
    @('522×8=4176', '739×7=5173'),
    @('759×6=4554', '509×8=4072'),
    @('999×6=5994', '562×8=4496'),
    @('946×2=1892', '495×5=2475'),
    @('663×9=5967', '364×6=2184'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
